$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.19%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.924"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.63%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07341"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.57%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.229"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "25.09%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.724"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.737"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.61%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9061"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.96%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09182"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "18.95%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1688"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.50%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08298"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.54%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03125"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.65%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09957"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.62%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005706"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.519"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.48%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.062"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.04%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3334"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.51%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.10%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.171"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.04%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2105"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.58%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04538"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.76%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.47%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004150"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.48%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.00%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003401"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-95.46%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.13%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04458"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.07%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007362"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.08%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009448"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "23.28%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1327"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002232"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.34%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009133"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.96%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006132"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.85%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.287"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.77%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002004"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-33.27%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
